$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7264
$ws1.Range("F3").Value = 406
$ws1.Range("F4").Value = 127
$ws1.Range("F5").Value = 179
$ws1.Range("F7").Value = 97
$ws1.Range("F8").Value = 615

# Sheet "全部类型" (all types) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7264
$ws4.Range("F3").Value = 406
$ws4.Range("F5").Value = 127
$ws4.Range("F6").Value = 179
$ws4.Range("F9").Value = 97
$ws4.Range("F10").Value = 615
